# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country names for rows 82 (was Sudan) and 83 (was Madagascar) ---
# After the edit, row 82 = Madagascar, row 83 = Sudan (values follow below).
$ws.Range("A82").Value = "Madagascar"
$ws.Range("A83").Value = "Sudan"

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 4 de Agosto de 2020 a las 12:50"

# --- Update numeric data per country row ---

# Row 4: Estados Unidos
$ws.Range("E4").Value = 2255250
$ws.Range("G4").Value = 39
$ws.Range("H4").Value = 158968

# Row 6: India
$ws.Range("B6").Value = 1861821
$ws.Range("C6").Value = 6490
$ws.Range("D6").Value = 1233589
$ws.Range("E6").Value = 589188
$ws.Range("G6").Value = 73
$ws.Range("H6").Value = 39044

# Row 14: Iran
$ws.Range("B14").Value = 314786
$ws.Range("C14").Value = 2751
$ws.Range("D14").Value = 272535
$ws.Range("E14").Value = 24634
$ws.Range("G14").Value = 212
$ws.Range("H14").Value = 17617

# Row 43: Emiratos Arabes Unidos
$ws.Range("B43").Value = 61352
$ws.Range("C43").Value = 189
$ws.Range("D43").Value = 55090
$ws.Range("E43").Value = 5911

# Row 45: Rumania
$ws.Range("B45").Value = 55241
$ws.Range("C45").Value = 1232
$ws.Range("D45").Value = 28006
$ws.Range("E45").Value = 24755
$ws.Range("G45").Value = 48
$ws.Range("H45").Value = 2480

# Row 58: Suiza
$ws.Range("B58").Value = 35746
$ws.Range("C58").Value = 130
$ws.Range("E58").Value = 2265

# Row 79: Estado de Palestina
$ws.Range("B79").Value = 12770
$ws.Range("C79").Value = 229
$ws.Range("D79").Value = 6419
$ws.Range("E79").Value = 6267

# Row 82: Madagascar (new data)
$ws.Range("B82").Value = 11895
$ws.Range("C82").Value = 235
$ws.Range("D82").Value = 9286
$ws.Range("E82").Value = 2486
$ws.Range("G82").Value = 5
$ws.Range("H82").Value = 123

# Row 83: Sudan (takes on former Sudan numbers)
$ws.Range("B83").Value = 11780
$ws.Range("C83").Value = 42
$ws.Range("D83").Value = 6194
$ws.Range("E83").Value = 4823
$ws.Range("G83").Value = 11
$ws.Range("H83").Value = 763

# Row 85: Senegal
$ws.Range("B85").Value = 10432
$ws.Range("C85").Value = 46
$ws.Range("D85").Value = 6920
$ws.Range("E85").Value = 3298
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 214

# Row 153: Malta
$ws.Range("B153").Value = 890
$ws.Range("C153").Value = 16
$ws.Range("E153").Value = 215

# Row 179: Gibraltar
$ws.Range("B179").Value = 189
$ws.Range("C179").Value = 1
$ws.Range("E179").Value = 6
